$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows: at position 2, then at position 4 (post first insert)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(4).Insert()

# Copy hyperlink-cell formatting (style) from an already-shifted F cell (F5, originally F3)
# onto the two new F cells (F2 and F4) so they use the existing "Hyperlink" cell style
# instead of the engine auto-creating a duplicate style entry.
$hlStyle = $ws.Range("F5").Style()
$ws.Range("F2").Style = $hlStyle
$ws.Range("F4").Style = $hlStyle

# Row 2: brand-new record
$ws.Range("A2").Value = '2025-12-22 18:27:59'
$ws.Range("B2").Value = '建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5434128'
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = '🔥AI,Ai ◆開発'

# Row 4: brand-new record
$ws.Range("A4").Value = '2025-12-22 18:27:59'
$ws.Range("B4").Value = '【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5439158'
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = '🔥AI,Ai'

# Rows that were shifted down only need column A (timestamp) refreshed
$ws.Range("A3").Value = '2025-12-22 18:27:59'
$ws.Range("A5").Value = '2025-12-22 18:27:59'
$ws.Range("A6").Value = '2025-12-22 18:27:59'
$ws.Range("A7").Value = '2025-12-22 18:27:59'
$ws.Range("A8").Value = '2025-12-22 18:27:59'
$ws.Range("A9").Value = '2025-12-22 18:27:59'
$ws.Range("A10").Value = '2025-12-22 18:27:59'
$ws.Range("A11").Value = '2025-12-22 18:27:59'
$ws.Range("A12").Value = '2025-12-22 18:27:59'
$ws.Range("A13").Value = '2025-12-22 18:27:59'

# Rebuild hyperlinks from scratch to exactly match the final row order/targets
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5434128') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5427956') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5439158') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5217096') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5458419') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5458381') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5456658') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5454210') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5458447') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5458919') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5418064') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5341051') | Out-Null
